$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.526.12"
$ws.Range("E2").Value = "  +0.19%  "
$ws.Range("D3").Value = "2.471.04"
$ws.Range("E3").Value = "  -0.70%  "
$ws.Range("E4").Value = "  -0.11%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "314.89"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.45%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "92.04"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -3.02%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.549"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("E8").Value = "  -0.14%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.514"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +2.76%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "32.25"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -4.05%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0790"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +1.04%  "
$ws.Range("E12").Value = "  +0.82%  "
$ws.Range("D13").Value = "2.849.92"
$ws.Range("E13").Value = "  -0.80%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "6.86"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -2.11%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "15.99"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +3.22%  "
$ws.Range("D16").Value = "2.466.56"
$ws.Range("E16").Value = "  +0.44%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.770"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -3.04%  "
$ws.Range("D18").Value = "41.512.40"
$ws.Range("E18").Value = "  +0.15%  "
$ws.Range("E19").Value = "  +2.24%  "
$ws.Range("D20").Value = "0.0₃0948"
$ws.Range("E20").Value = "  +2.51%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "71.33"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +3.27%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "11.14"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -0.93%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "236.17"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -0.73%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "2.73"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -1.36%  "
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("E26").Value = "  -1.08%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "24.62"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +1.29%  "
$ws.Range("E28").Value = "  -0.69%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "9.68"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -0.57%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "35.41"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -3.08%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "156.32"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +2.52%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "5.45"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -0.84%  "
$ws.Range("E33").Value = "  -0.96%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.0759"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +1.13%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "17.31"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -2.28%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "2.89"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -6.20%  "
$ws.Range("E37").Value = "  +1.63%  "
$ws.Range("E38").Value = "  -0.57%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "1.80"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -4.42%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "2.25"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -12.04%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "4.06"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -4.02%  "
$ws.Range("E42").Value = "  -0.23%  "
$ws.Range("D43").Value = "1.944.44"
$ws.Range("E43").Value = "  -3.06%  "
$ws.Range("E44").Value = "  -1.43%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "18.48"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -5.25%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "2.94"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -3.01%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "9.08"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +3.04%  "
$ws.Range("D48").Value = "2.705.97"
$ws.Range("E48").Value = "  -1.07%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "97.13"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -0.62%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "67.29"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -4.01%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "52.43"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +2.76%  "
